# Insert a new data row at row 1071 (above the existing "Valencia / Primera"
# entry dated 44610), shifting all subsequent rows (old 1071-1152) down by
# one so they become rows 1072-1153. This grows the used range from
# A1:T1152 to A1:T1153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 1071 (and everything below it) down by one row.
$ws.Rows.Item(1071).Insert()

# Populate the newly inserted row 1071 with the new price-sheet entry.
$ws.Range("A1071").Value = 9
$ws.Range("B1071").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1071").Value = "Metropolitana"
$ws.Range("D1071").Value = 45013
$ws.Range("E1071").Value = 13
$ws.Range("F1071").Value = "Fruta"
$ws.Range("G1071").Value = 100102
$ws.Range("H1071").Value = "Cítricos"
$ws.Range("I1071").Value = 100102005
$ws.Range("J1071").Value = "Naranja"
$ws.Range("K1071").Value = "Valencia"
$ws.Range("L1071").Value = "Primera"
$ws.Range("M1071").Value = 400
$ws.Range("N1071").Value = 12000
$ws.Range("O1071").Value = 12000
$ws.Range("P1071").Value = 12000
$ws.Range("Q1071").Value = "$/caja 15 kilos granel"
$ws.Range("R1071").Value = "Región de O'Higgins"
$ws.Range("S1071").Value = 800
$ws.Range("T1071").Value = 15
